$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new day's price (11-02-2026) was published. It is identical in every
# respect (description, product code, basic price, circular date/link) to
# what was previously the most recent row, so the simplest, most faithful
# way to reproduce the diff is to duplicate the current top data row (row 2)
# into a freshly inserted row above it - this naturally pushes every
# existing row down by one (matching the diff's "date shifted down" pattern)
# while preserving cell styles, number formats and the existing hyperlink.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(2).Insert()

# The engine sometimes forgets the "Basic Price" number format (0.000) on a
# freshly inserted row, so make sure it still matches the rest of column D.
$ws.Range("D2").NumberFormat = "0.000"

# Only the Date column actually differs for the brand-new top row; force it
# to stay plain text (matching every other Date cell) instead of letting it
# be auto-converted into a date serial value.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "11-02-2026"

# Because row 189 was the last row in the sheet, shifting everything down by
# one means the former row 189 now also has to exist at row 190 (duplicate
# of the oldest circular, 07-08-2025) to keep 189 data rows total. Row 190's
# values/format already come through from the cascade above; it just needs
# its own hyperlink on column F, since hyperlinks aren't auto-extended past
# the original last row.
$lastUrl = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Hyperlinks.Add($ws.Range("F190"), $lastUrl)

# Re-apply F189's formatting to F190 so the newly-added hyperlink doesn't
# leave behind Excel's default blue/underlined "Hyperlink" style, keeping it
# visually consistent with every other Circular Link cell in the sheet.
$ws.Range("F189").Copy()
$ws.Range("F190").PasteSpecial(-4122)
